$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimates")

# Mark these stories as "Completed" (checking the Completed column in Table1) -
# the calculated "Completed Points"/"Completed Hours" columns recompute automatically.
$completedRows = 24,25,27,28,52,53,54,57,62,71,72,73,74
foreach ($r in $completedRows) {
    $ws.Range("E$r").Value = $true
}

# Tidy up stray cell formatting on the "Completed Points"/"Completed Hours" columns
# for the not-yet-estimated backlog rows so they match the rest of the column.
$ws.Range("F88:G95").ClearFormats()

# Leave the selection on the last checkbox that was toggled.
$ws.Activate() | Out-Null
$ws.Range("E57").Select() | Out-Null
